$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4214.4585
$ws.Range("I70").Value = 2001.3334
$ws.Range("J70").Value = 4530.619
$ws.Range("K70").Value = 6004.0002
$ws.Range("L70").Value = 13591.857
$ws.Range("M70").Value = -5734.0002
$ws.Range("N70").Value = -14131.857
$ws.Range("H73").Value = 4214.4585
$ws.Range("I73").Value = 2001.3334
$ws.Range("J73").Value = 4530.619
$ws.Range("K73").Value = 6004.0002
$ws.Range("L73").Value = 13591.857
$ws.Range("M73").Value = -5068.0002
$ws.Range("N73").Value = -15463.857
$ws.Range("H132").Value = 1626.9429
$ws.Range("I132").Value = 1457.5938
$ws.Range("J132").Value = 3433.3333
$ws.Range("K132").Value = 4372.7814
$ws.Range("L132").Value = 10299.9999
$ws.Range("M132").Value = -1842.7814
$ws.Range("N132").Value = -15359.9999
$ws.Range("H138").Value = 2930.5278
$ws.Range("J138").Value = 2782.6667
$ws.Range("L138").Value = 8348.000100000001
$ws.Range("N138").Value = -18628.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5249.9
$ws.Range("I45").Value = 3750
$ws.Range("K45").Value = 3750
$ws.Range("M45").Value = -3373
$ws.Range("H74").Value = 6108.1934
$ws.Range("I74").Value = 4999.625
$ws.Range("K74").Value = 4999.625
$ws.Range("M74").Value = -4125.625
$ws.Range("H77").Value = 6108.1934
$ws.Range("I77").Value = 4999.625
$ws.Range("K77").Value = 24998.125
$ws.Range("M77").Value = -20630.125
$ws.Range("H110").Value = 4504.645
$ws.Range("I110").Value = 3195.3333
$ws.Range("J110").Value = 8993.714
$ws.Range("K110").Value = 3195.3333
$ws.Range("L110").Value = 8993.714
$ws.Range("M110").Value = -1150.3333
$ws.Range("N110").Value = -13083.714

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2118.4243
$ws.Range("I86").Value = 1889.7693
$ws.Range("K86").Value = 1889.7693
$ws.Range("M86").Value = -766.7692999999999
$ws.Range("H89").Value = 2118.4243
$ws.Range("I89").Value = 1889.7693
$ws.Range("K89").Value = 9448.8465
$ws.Range("M89").Value = -3832.8465
$ws.Range("H99").Value = 4099.4585
$ws.Range("I99").Value = 3080.6924
$ws.Range("J99").Value = 5303.4546
$ws.Range("K99").Value = 3080.6924
$ws.Range("L99").Value = 5303.4546
$ws.Range("M99").Value = -1582.6924
$ws.Range("N99").Value = -8299.454600000001
$ws.Range("H107").Value = 3272.1428
$ws.Range("I107").Value = 3098.4
$ws.Range("K107").Value = 3098.4
$ws.Range("M107").Value = -1178.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3523.7896
$ws.Range("I16").Value = 1997.3572
$ws.Range("K16").Value = 1997.3572
$ws.Range("M16").Value = -1710.3572
$ws.Range("H31").Value = 3770.6538
$ws.Range("I31").Value = 2819.7144
$ws.Range("J31").Value = 4880.0835
$ws.Range("K31").Value = 2819.7144
$ws.Range("L31").Value = 4880.0835
$ws.Range("M31").Value = -2524.7144
$ws.Range("N31").Value = -5470.0835
$ws.Range("H34").Value = 3770.6538
$ws.Range("I34").Value = 2819.7144
$ws.Range("J34").Value = 4880.0835
$ws.Range("K34").Value = 2819.7144
$ws.Range("L34").Value = 4880.0835
$ws.Range("M34").Value = -2617.7144
$ws.Range("N34").Value = -5284.0835
$ws.Range("H99").Value = 5117.25
$ws.Range("I99").Value = 5499.6665
$ws.Range("K99").Value = 5499.6665
$ws.Range("M99").Value = -4001.6665
$ws.Range("H113").Value = 3523.7896
$ws.Range("I113").Value = 1997.3572
$ws.Range("K113").Value = 1997.3572
$ws.Range("M113").Value = 172.6428000000001
$ws.Range("H126").Value = 5117.25
$ws.Range("I126").Value = 5499.6665
$ws.Range("K126").Value = 16498.9995
$ws.Range("M126").Value = -14028.9995
$ws.Range("H140").Value = 179998
$ws.Range("J140").Value = 179998
$ws.Range("L140").Value = 179998
$ws.Range("N140").Value = -190358

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1565.5385
$ws.Range("J5").Value = 1916.6666
$ws.Range("L5").Value = 5749.9998
$ws.Range("N5").Value = -5973.9998
$ws.Range("H55").Value = 813.25
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H92").Value = 245.4
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 175.66667
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 527.00001
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -3023.00001
$ws.Range("H126").Value = 13833.167
$ws.Range("I126").Value = 12666.333
$ws.Range("K126").Value = 37998.999
$ws.Range("M126").Value = -33058.999
$ws.Range("H135").Value = 1565.5385
$ws.Range("J135").Value = 1916.6666
$ws.Range("L135").Value = 17249.9994
$ws.Range("N135").Value = -22319.9994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5713.364
$ws.Range("I122").Value = 5513.3335
$ws.Range("J122").Value = 5953.4
$ws.Range("K122").Value = 16540.0005
$ws.Range("L122").Value = 17860.2
$ws.Range("M122").Value = -14090.0005
$ws.Range("N122").Value = -22760.2
$ws.Range("H126").Value = 5999.8335
$ws.Range("I126").Value = 5166.6665
$ws.Range("J126").Value = 6833
$ws.Range("K126").Value = 15499.9995
$ws.Range("L126").Value = 20499
$ws.Range("M126").Value = -13029.9995
$ws.Range("N126").Value = -25439
$ws.Range("H132").Value = 1068.2858
$ws.Range("I132").Value = 1079.6666
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3238.9998
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -708.9998000000001
$ws.Range("N132").Value = -8060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2403.9092
$ws.Range("I7").Value = 2403.9092
$ws.Range("K7").Value = 2403.9092
$ws.Range("M7").Value = -2291.9092
$ws.Range("H16").Value = 1739.9474
$ws.Range("I16").Value = 1603.6875
$ws.Range("K16").Value = 1603.6875
$ws.Range("M16").Value = -1433.6875
$ws.Range("H61").Value = 144399.86
$ws.Range("I61").Value = 251699.75
$ws.Range("K61").Value = 251699.75
$ws.Range("M61").Value = -251497.75
$ws.Range("H113").Value = 144399.86
$ws.Range("I113").Value = 251699.75
$ws.Range("K113").Value = 251699.75
$ws.Range("M113").Value = -249529.75
$ws.Range("H126").Value = 2403.9092
$ws.Range("I126").Value = 2403.9092
$ws.Range("K126").Value = 7211.7276
$ws.Range("M126").Value = -4741.7276
$ws.Range("H132").Value = 7525.92
$ws.Range("I132").Value = 7978.7803
$ws.Range("J132").Value = 5462.8887
$ws.Range("K132").Value = 23936.3409
$ws.Range("L132").Value = 16388.6661
$ws.Range("M132").Value = -21406.3409
$ws.Range("N132").Value = -21448.6661

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5741.0625
$ws.Range("I62").Value = 5205.5
$ws.Range("J62").Value = 9490
$ws.Range("K62").Value = 5205.5
$ws.Range("L62").Value = 9490
$ws.Range("M62").Value = -4581.5
$ws.Range("N62").Value = -10738
$ws.Range("H65").Value = 5741.0625
$ws.Range("I65").Value = 5205.5
$ws.Range("J65").Value = 9490
$ws.Range("K65").Value = 26027.5
$ws.Range("L65").Value = 47450
$ws.Range("M65").Value = -22907.5
$ws.Range("N65").Value = -53690
$ws.Range("H122").Value = 4182.485
$ws.Range("I122").Value = 3911.16
$ws.Range("K122").Value = 11733.48
$ws.Range("M122").Value = -9283.48
$ws.Range("H126").Value = 3135.92
$ws.Range("I126").Value = 2976.2222
$ws.Range("K126").Value = 8928.6666
$ws.Range("M126").Value = -6458.6666
$ws.Range("H132").Value = 3859.75
$ws.Range("I132").Value = 4037.4333
$ws.Range("J132").Value = 3326.7
$ws.Range("K132").Value = 12112.2999
$ws.Range("L132").Value = 9980.099999999999
$ws.Range("M132").Value = -9582.2999
$ws.Range("N132").Value = -15040.1
